$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh rotates the daily price records across the dated rows:
# row 2 <- old row 6, row 3 <- old row 5, row 5 <- old row 2, row 6 <- old row 3
# (row 4 keeps its values). Capture the "before" values first so the
# in-place writes below don't clobber data we still need to read.

$cols = @("D","J","K","L","M","P")

$before2 = @{}
$before3 = @{}
$before5 = @{}
$before6 = @{}

foreach ($col in $cols) {
    $before2[$col] = $ws.Range("$col" + "2").Value()
    $before3[$col] = $ws.Range("$col" + "3").Value()
    $before5[$col] = $ws.Range("$col" + "5").Value()
    $before6[$col] = $ws.Range("$col" + "6").Value()
}

foreach ($col in $cols) {
    $ws.Range("$col" + "2").Value = $before6[$col]
    $ws.Range("$col" + "3").Value = $before5[$col]
    $ws.Range("$col" + "5").Value = $before2[$col]
    $ws.Range("$col" + "6").Value = $before3[$col]
}
